$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 95, pushing the existing rows 95-109 down to 97-111.
# (A new week of "Chirimoya" price records was added at the top of this market's block.)
$ws.Rows.Item(95).Insert()
$ws.Rows.Item(95).Insert()

# --- New row 95 ---
$ws.Cells.Item(95, 1).Value = 10
$ws.Cells.Item(95, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(95, 3).Value = "La Araucanía"
$ws.Cells.Item(95, 4).Value = 44522
$ws.Cells.Item(95, 5).Value = 9
$ws.Cells.Item(95, 6).Value = "Fruta"
$ws.Cells.Item(95, 7).Value = 100107
$ws.Cells.Item(95, 8).Value = "Otros"
$ws.Cells.Item(95, 9).Value = 100107002
$ws.Cells.Item(95, 10).Value = "Chirimoya"
$ws.Cells.Item(95, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(95, 12).Value = "Primera"
$ws.Cells.Item(95, 13).Value = 50
$ws.Cells.Item(95, 14).Value = 3000
$ws.Cells.Item(95, 15).Value = 3000
$ws.Cells.Item(95, 16).Value = 3000
$ws.Cells.Item(95, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(95, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(95, 19).Value = 3000
$ws.Cells.Item(95, 20).Value = 1

# --- New row 96 ---
$ws.Cells.Item(96, 1).Value = 10
$ws.Cells.Item(96, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(96, 3).Value = "La Araucanía"
$ws.Cells.Item(96, 4).Value = 44522
$ws.Cells.Item(96, 5).Value = 9
$ws.Cells.Item(96, 6).Value = "Fruta"
$ws.Cells.Item(96, 7).Value = 100107
$ws.Cells.Item(96, 8).Value = "Otros"
$ws.Cells.Item(96, 9).Value = 100107002
$ws.Cells.Item(96, 10).Value = "Chirimoya"
$ws.Cells.Item(96, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(96, 12).Value = "Segunda"
$ws.Cells.Item(96, 13).Value = 20
$ws.Cells.Item(96, 14).Value = 2500
$ws.Cells.Item(96, 15).Value = 2500
$ws.Cells.Item(96, 16).Value = 2500
$ws.Cells.Item(96, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(96, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(96, 19).Value = 2500
$ws.Cells.Item(96, 20).Value = 1
